$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of license data for Vansh Tyagi
$ws.Range("A5").Value = "Vansh Tyagi"
$ws.Range("B5").Value = "UP14 20230028483"
